$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "44-39=5"
$t.Cell(1,2).Range.Text = "63-53=10"
$t.Cell(1,3).Range.Text = "41-2=39"
$t.Cell(1,4).Range.Text = "99-39=60"
$t.Cell(1,5).Range.Text = "2+18=20"
$t.Cell(2,1).Range.Text = "56-24=32"
$t.Cell(2,2).Range.Text = "57-46=11"
$t.Cell(2,3).Range.Text = "64+31=95"
$t.Cell(2,4).Range.Text = "1+28=29"
$t.Cell(2,5).Range.Text = "45+49=94"
$t.Cell(3,1).Range.Text = "66-3=63"
$t.Cell(3,2).Range.Text = "36+39=75"
$t.Cell(3,3).Range.Text = "27-19=8"
$t.Cell(3,4).Range.Text = "11+72=83"
$t.Cell(3,5).Range.Text = "40-14=26"
$t.Cell(4,1).Range.Text = "48-13=35"
$t.Cell(4,2).Range.Text = "73+18=91"
$t.Cell(4,3).Range.Text = "21+15=36"
$t.Cell(4,4).Range.Text = "92-70=22"
$t.Cell(4,5).Range.Text = "35+7=42"
$t.Cell(5,1).Range.Text = "47+28=75"
$t.Cell(5,2).Range.Text = "43-29=14"
$t.Cell(5,3).Range.Text = "91-49=42"
$t.Cell(5,4).Range.Text = "29+68=97"
$t.Cell(5,5).Range.Text = "80-49=31"
$t.Cell(6,1).Range.Text = "28+27=55"
$t.Cell(6,2).Range.Text = "75-54=21"
$t.Cell(6,3).Range.Text = "6+91=97"
$t.Cell(6,4).Range.Text = "6+64=70"
$t.Cell(6,5).Range.Text = "91-54=37"
$t.Cell(7,1).Range.Text = "62-1=61"
$t.Cell(7,2).Range.Text = "44+39=83"
$t.Cell(7,3).Range.Text = "29+69=98"
$t.Cell(7,4).Range.Text = "9+65=74"
$t.Cell(7,5).Range.Text = "73-59=14"
$t.Cell(8,1).Range.Text = "73-10=63"
$t.Cell(8,2).Range.Text = "51-50=1"
$t.Cell(8,3).Range.Text = "70-21=49"
$t.Cell(8,4).Range.Text = "48+34=82"
$t.Cell(8,5).Range.Text = "96-93=3"
$t.Cell(9,1).Range.Text = "13+41=54"
$t.Cell(9,2).Range.Text = "48-28=20"
$t.Cell(9,3).Range.Text = "12+2=14"
$t.Cell(9,4).Range.Text = "0+98=98"
$t.Cell(9,5).Range.Text = "75-22=53"
$t.Cell(10,1).Range.Text = "86-74=12"
$t.Cell(10,2).Range.Text = "87-27=60"
$t.Cell(10,3).Range.Text = "27+3=30"
$t.Cell(10,4).Range.Text = "83-43=40"
$t.Cell(10,5).Range.Text = "27+29=56"
$t.Cell(11,1).Range.Text = "15+21=36"
$t.Cell(11,2).Range.Text = "34+2=36"
$t.Cell(11,3).Range.Text = "34+55=89"
$t.Cell(11,4).Range.Text = "51-48=3"
$t.Cell(11,5).Range.Text = "77-14=63"
$t.Cell(12,1).Range.Text = "55+28=83"
$t.Cell(12,2).Range.Text = "57+32=89"
$t.Cell(12,3).Range.Text = "17+25=42"
$t.Cell(12,4).Range.Text = "52+1=53"
$t.Cell(12,5).Range.Text = "2+38=40"
$t.Cell(13,1).Range.Text = "8+32=40"
$t.Cell(13,2).Range.Text = "58-26=32"
$t.Cell(13,3).Range.Text = "18+69=87"
$t.Cell(13,4).Range.Text = "2+45=47"
$t.Cell(13,5).Range.Text = "80-34=46"
$t.Cell(14,1).Range.Text = "40+49=89"
$t.Cell(14,2).Range.Text = "53-41=12"
$t.Cell(14,3).Range.Text = "46+49=95"
$t.Cell(14,4).Range.Text = "12+56=68"
$t.Cell(14,5).Range.Text = "20-11=9"
$t.Cell(15,1).Range.Text = "25-17=8"
$t.Cell(15,2).Range.Text = "13+32=45"
$t.Cell(15,3).Range.Text = "89-41=48"
$t.Cell(15,4).Range.Text = "83+15=98"
$t.Cell(15,5).Range.Text = "45-19=26"
$t.Cell(16,1).Range.Text = "47+36=83"
$t.Cell(16,2).Range.Text = "51+4=55"
$t.Cell(16,3).Range.Text = "10+47=57"
$t.Cell(16,4).Range.Text = "25+23=48"
$t.Cell(16,5).Range.Text = "8+76=84"
$t.Cell(17,1).Range.Text = "46+6=52"
$t.Cell(17,2).Range.Text = "12+36=48"
$t.Cell(17,3).Range.Text = "1+54=55"
$t.Cell(17,4).Range.Text = "55+26=81"
$t.Cell(17,5).Range.Text = "15+54=69"
$t.Cell(18,1).Range.Text = "21+3=24"
$t.Cell(18,2).Range.Text = "63+34=97"
$t.Cell(18,3).Range.Text = "95-11=84"
$t.Cell(18,4).Range.Text = "56-16=40"
$t.Cell(18,5).Range.Text = "26-8=18"
$t.Cell(19,1).Range.Text = "73-39=34"
$t.Cell(19,2).Range.Text = "98-63=35"
$t.Cell(19,3).Range.Text = "89-51=38"
$t.Cell(19,4).Range.Text = "3-2=1"
$t.Cell(19,5).Range.Text = "56+9=65"
$t.Cell(20,1).Range.Text = "15+34=49"
$t.Cell(20,2).Range.Text = "98-62=36"
$t.Cell(20,3).Range.Text = "10+40=50"
$t.Cell(20,4).Range.Text = "40-13=27"
$t.Cell(20,5).Range.Text = "34+40=74"
